# Week 15 simulations update
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet: append new simulation-week numbers to the four long number lists
# ---------------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Text + " 3 -1 11 8 -1 12 8 4 2 4 1 4 1 4 -5 -1 5 11 17 4 0 1 2 -2 2 4"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Text + " 5 -1 0 2 6 5 4 11 0 4 -1 4 11 4 0 11 1 4 3 11 8 5 9 0 2 1 -2 1 11 11 4 2 3 5 3"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Text + " 7 11 15 22 10 23 9 8 14 6 18 7 9 3 19 19 13 18 11"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Text + " -2 6 36 9 13 3 20 4 4 14 12 7 14 3 11 6 2 23"

# ---------------------------------------------------------------------------
# OFF sheet: updated season totals (row 2 = Home, row 3 = Road)
# ---------------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("C2").Value = 377
$offWs.Range("E2").Value = 12
$offWs.Range("F2").Value = 87
$offWs.Range("G2").Value = 109
$offWs.Range("I2").Value = 15
$offWs.Range("J2").Value = 60
$offWs.Range("L2").Value = 510
$offWs.Range("M2").Value = 314
$offWs.Range("O2").Value = 50
$offWs.Range("P2").Value = 26
$offWs.Range("Q2").Value = 921

$offWs.Range("B3").Value = 20
$offWs.Range("C3").Value = 299
$offWs.Range("E3").Value = 68
$offWs.Range("F3").Value = 205
$offWs.Range("I3").Value = 112
$offWs.Range("J3").Value = 80
$offWs.Range("N3").Value = 31

# ---------------------------------------------------------------------------
# DEF sheet: updated season totals (row 2 = Home, row 3 = Road)
# ---------------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("B2").Value = 9
$defWs.Range("C2").Value = 361
$defWs.Range("E2").Value = 16
$defWs.Range("F2").Value = 97
$defWs.Range("G2").Value = 97
$defWs.Range("J2").Value = 54
$defWs.Range("L2").Value = 472
$defWs.Range("M2").Value = 314
$defWs.Range("O2").Value = 29
$defWs.Range("P2").Value = 16
$defWs.Range("Q2").Value = 918

$defWs.Range("B3").Value = 10
$defWs.Range("C3").Value = 286
$defWs.Range("D3").Value = 11
$defWs.Range("E3").Value = 63
$defWs.Range("F3").Value = 206
$defWs.Range("H3").Value = 53
$defWs.Range("I3").Value = 93

# ---------------------------------------------------------------------------
# ST sheet: updated totals + appended simulation-week numbers
# ---------------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value = 139
$stWs.Range("D2").Value = 97
$stWs.Range("F2").Value = 134
$stWs.Range("G2").Value = 127
$stWs.Range("B3").Value = 111

$stWs.Range("D3").Value = $stWs.Range("D3").Text + " 54 33"
$stWs.Range("B4").Value = $stWs.Range("B4").Text + " 63"
$stWs.Range("D4").Value = $stWs.Range("D4").Text + " 0 0"
$stWs.Range("B5").Value = $stWs.Range("B5").Text + " 38"
$stWs.Range("D5").Value = $stWs.Range("D5").Text + " 0 0"
$stWs.Range("B6").Value = $stWs.Range("B6").Text + " 26 28 0 25"

# ---------------------------------------------------------------------------
# TURNS sheet: updated totals
# ---------------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")

$turnsWs.Range("B2").Value = 22
$turnsWs.Range("D2").Value = 12
$turnsWs.Range("E2").Value = 17
$turnsWs.Range("E3").Value = 16

# ---------------------------------------------------------------------------
# PEN sheet: updated totals
# ---------------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")

$penWs.Range("B2").Value = 30
